$wb = $excel.ActiveWorkbook

# Metadata sheet is the first sheet in the workbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" property value (row 8, column B)
$ws.Range("B8").Value = "2025-01-24T13:24:43+00:00"

# Fill in the previously-empty "Description" property value (row 12, column B)
$ws.Range("B12").Value = "ValueSet regroupant des valuesets du NOS pour le code de la division territoriale"
